$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G20").Value = 16
$ws.Range("G21").Value = 5
$ws.Range("G27").Value = 10
$ws.Range("G28").Value = 5
$ws.Range("G37").Value = 16
$ws.Range("G38").Value = 5
$ws.Range("G39").Value = 3
$ws.Range("G40").Value = 19
$ws.Range("G41").Value = 10
$ws.Range("G42").Value = 5
$ws.Range("G43").Value = 12
$ws.Range("G44").Value = 14

$ws.StandardWidth = 11.58984375

$ws.Range("G45").Select()
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
